$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Key Features section updates ---

# Row 14: "Support for Hot-Swapping multiple levels without a restart"
# previously had a stray "y" mark in column E -> remove it entirely
$ws.Range("E14").Clear()

# Row 15: "Support for adding Point and/or Spotlight sources via Blender"
# previously had a stray "?" mark in column E -> remove it entirely
$ws.Range("E15").Clear()

# Row 19: "Support for Rendering Bounding Volume/Collission Data" (Key tally)
# count of Key features picked goes from 22 -> 33
$ws.Range("G19").Value = 33

# Row 22: "Optimizing your Renderer with Frustum Culling"
# mark this Key feature as done ("x" in column D, matching the style used
# by the other completed rows), and remove the stray "?" mark in column E
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D22").Value = "x"
$ws.Range("E22").Clear()

# Row 24: "Support for Diffuse and Specular Textured Materials in your Level"
# previously had a stray "?" mark in column E -> remove it entirely
$ws.Range("E24").Clear()

# Row 46: "Play some Level Music and add at least one Sound Effect (GMusic & GSound)"
# previously had a stray "?" mark in column E -> remove it entirely
$ws.Range("E46").Clear()

$excel.CutCopyMode = 0

# --- View state updates ---
$ws.Range("F24").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 10
